# Refresh the crypto market snapshot on Sheet1: latest Price (column D) and
# 1h Volume change (column E) for every coin row. A handful of rows also swap
# which coin occupies that rank (Coin name in B + Link in C) to reflect the
# reshuffled leaderboard order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.321.25"
$ws.Range("E2").Value = "  +2.65%  "

# Row 3
$ws.Range("D3").Value = "1.915.89"
$ws.Range("E3").Value = "  +1.32%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "'248.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "

# Row 6
$ws.Range("D6").Value = "'0.694"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "

# Row 8
$ws.Range("D8").Value = "'43.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.58%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'57.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.21%  "

# Row 10
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "'0.364"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.01%  "

# Row 11
$ws.Range("D11").Value = "'0.0765"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.14%  "

# Row 12
$ws.Range("D12").Value = "'0.0993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.51%  "

# Row 13
$ws.Range("D13").Value = "'14.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.69%  "

# Row 14
$ws.Range("D14").Value = "'0.796"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.02%  "

# Row 15
$ws.Range("D15").Value = "2.183.59"
$ws.Range("E15").Value = "  +0.85%  "

# Row 16
$ws.Range("D16").Value = "'5.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.08%  "

# Row 17
$ws.Range("D17").Value = "1.905.59"
$ws.Range("E17").Value = "  +0.69%  "

# Row 18
$ws.Range("D18").Value = "36.207.80"
$ws.Range("E18").Value = "  +2.26%  "

# Row 19
$ws.Range("D19").Value = "'74.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.26%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0842"
$ws.Range("E20").Value = "  +2.27%  "

# Row 21
$ws.Range("D21").Value = "'250.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.21%  "

# Row 22
$ws.Range("D22").Value = "'13.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.13%  "

# Row 23
$ws.Range("D23").Value = "'5.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.98%  "

# Row 24
$ws.Range("E24").Value = "  +1.17%  "

# Row 25
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("D26").Value = "'2.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.92%  "

# Row 27
$ws.Range("D27").Value = "'167.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "

# Row 28
$ws.Range("D28").Value = "'8.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.72%  "

# Row 29
$ws.Range("D29").Value = "'18.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.07%  "

# Row 30
$ws.Range("D30").Value = "'0.129"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "

# Row 31
$ws.Range("D31").Value = "'4.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.09%  "

# Row 32
$ws.Range("D32").Value = "'0.0611"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.13%  "

# Row 33
$ws.Range("D33").Value = "'4.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.22%  "

# Row 34
$ws.Range("D34").Value = "'1.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.82%  "

# Row 35
$ws.Range("E35").Value = "  -0.04%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.0848"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +22.30%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -16.59%  "

# Row 38
$ws.Range("D38").Value = "'0.859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.33%  "

# Row 39
$ws.Range("D39").Value = "'2.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.87%  "

# Row 40
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'103.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.91%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0229"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.18%  "

# Row 42
$ws.Range("D42").Value = "'15.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.72%  "

# Row 43
$ws.Range("D43").Value = "'17.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.80%  "

# Row 44
$ws.Range("D44").Value = "'1.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.96%  "

# Row 45
$ws.Range("D45").Value = "1.342.23"
$ws.Range("E45").Value = "  +3.76%  "

# Row 46
$ws.Range("D46").Value = "'2.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.81%  "

# Row 47
$ws.Range("D47").Value = "'0.0809"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.67%  "

# Row 48
$ws.Range("E48").Value = "  +0.72%  "

# Row 49
$ws.Range("E49").Value = "  +1.58%  "

# Row 50
$ws.Range("D50").Value = "'6.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.42%  "

# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'42.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.89%  "
